# Fix "Año" (E), "Mes" (F) and "Año.Mes" (G) columns so that they are
# correctly derived from the "Fecha" column (I), which is stored as
# text in "DD/MM/YYYY" format. Previously the day portion was being
# used as the month by mistake.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 9).End(-4162).Row  # xlUp = -4162, column I = 9

for ($r = 2; $r -le $lastRow; $r++) {
    $fecha = $ws.Cells.Item($r, 9).Text
    if ([string]::IsNullOrEmpty($fecha)) {
        continue
    }

    $parts = $fecha.Split('/')
    if ($parts.Count -ne 3) {
        continue
    }

    $day = [int]$parts[0]
    $month = [int]$parts[1]
    $year = [int]$parts[2]

    $ws.Cells.Item($r, 5).Value = $year
    $ws.Cells.Item($r, 6).Value = $month

    if ($month -lt 10) {
        $monthPadded = "0$month"
    } else {
        $monthPadded = "$month"
    }

    $gCell = $ws.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = "$year.$monthPadded"
}
